$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in the title cell
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 20:57"

# Row 4 - Estados Unidos (rank unchanged)
$ws.Range("B4").Value = 3653378
$ws.Range("C4").Value = 36551
$ws.Range("D4").Value = 1661097
$ws.Range("E4").Value = 1851708
$ws.Range("G4").Value = 429
$ws.Range("H4").Value = 140573

# Row 6 - India (rank unchanged)
$ws.Range("B6").Value = 1004647
$ws.Range("C6").Value = 34478
$ws.Range("D6").Value = 636569
$ws.Range("E6").Value = 342469
$ws.Range("G6").Value = 680
$ws.Range("H6").Value = 25609

# Rows 9-11: Sudafrica moved above Chile/Mexico (re-sort by total cases)
# Row 9 now shows Sudafrica with updated figures
$ws.Range("A9").Value = "Sudafrica"
$ws.Range("B9").Value = 324221
$ws.Range("C9").Value = 13172
$ws.Range("D9").Value = 165591
$ws.Range("E9").Value = 153961
$ws.Range("G9").Value = 216
$ws.Range("H9").Value = 4669

# Row 10 now shows Chile (previous Chile figures)
$ws.Range("A10").Value = "Chile"
$ws.Range("B10").Value = 323698
$ws.Range("C10").Value = 2493
$ws.Range("D10").Value = 295301
$ws.Range("E10").Value = 21107
$ws.Range("G10").Value = 104
$ws.Range("H10").Value = 7290

# Row 11 now shows Mexico (previous Mexico figures)
$ws.Range("A11").Value = "Mexico"
$ws.Range("B11").Value = 317635
$ws.Range("C11").Value = 6149
$ws.Range("D11").Value = 199129
$ws.Range("E11").Value = 81600
$ws.Range("G11").Value = 579
$ws.Range("H11").Value = 36906

# Row 39 - Emiratos Arabes Unidos
$ws.Range("B39").Value = 56129
$ws.Range("C39").Value = 281
$ws.Range("D39").Value = 47412
$ws.Range("E39").Value = 8382

# Row 58 - Irlanda
$ws.Range("B58").Value = 25698
$ws.Range("C58").Value = 15
$ws.Range("E58").Value = 585
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1749

# Row 65 - Marruecos
$ws.Range("B65").Value = 16545
$ws.Range("C65").Value = 283
$ws.Range("D65").Value = 13965
$ws.Range("E65").Value = 2317
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 263

# Row 162 - Vietnam
$ws.Range("D162").Value = 356
$ws.Range("E162").Value = 25

# Shared-string reorder: Groenlandia moved above Islas Malvinas
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
